$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = '2017-12-24 17_02 - Dr. Amit'
$ws.Range("C4").Value = 2.32
$ws.Range("D4").Value = 60.28
$ws.Range("E4").Value = 97.41
$ws.Range("F4").Value = 92.95
$ws.Range("G4").Value = 2.61
$ws.Range("I4").Value = 81216
$ws.Range("J4").Value = 40.36
$ws.Range("K4").Value = 78406
$ws.Range("L4").Value = 39.28
$ws.Range("M4").Value = 82329
$ws.Range("N4").Value = 40.79
$ws.Range("O4").Value = 5.21
$ws.Range("P4").Value = 5.25
$ws.Range("Q4").Value = 5.07
$ws.Range("R4").Value = 4.51
$ws.Range("S4").Value = 5.71
$ws.Range("T4").Value = 5.01
$ws.Range("U4").Value = 4.99
$ws.Range("V4").Value = 92.31999999999999
$ws.Range("W4").Value = 79.12
$ws.Range("X4").Value = 93.67
$ws.Range("Y4").Value = 98.93000000000001
$ws.Range("Z4").Value = 92.47
$ws.Range("AA4").Value = 98.73999999999999
$ws.Range("AB4").Value = 96.09
$ws.Range("AC4").Value = 0.23
$ws.Range("AD4").Value = 0.63
$ws.Range("AE4").Value = -0.19
$ws.Range("AF4").Value = 0.03
$ws.Range("AG4").Value = -0.23
$ws.Range("AH4").Value = 0.04
$ws.Range("AI4").Value = 0.12
$ws.Range("AJ4").Value = 93.05
$ws.Range("AK4").Value = 79.20999999999999
$ws.Range("AL4").Value = 50.05
$ws.Range("AM4").Value = 49.95
$ws.Range("AN4").Value = 4.8
$ws.Range("AO4").Value = 4.69
$ws.Range("AP4").Value = 4.25
$ws.Range("AQ4").Value = 6.33
$ws.Range("AR4").Value = 4.75
$ws.Range("AS4").Value = 5.84
$ws.Range("AT4").Value = 4.74
$ws.Range("AU4").Value = 5.12
$ws.Range("AV4").Value = 5.23
$ws.Range("AW4").Value = 5.34
$ws.Range("AX4").Value = 5.35
$ws.Range("AY4").Value = 4.84
$ws.Range("AZ4").Value = 5.74
$ws.Range("BA4").Value = 5.79
$ws.Range("BB4").Value = 3.82
$ws.Range("BC4").Value = 4.85
$ws.Range("BD4").Value = 5.84
$ws.Range("BE4").Value = 4.8
$ws.Range("BF4").Value = 4.45
$ws.Range("BH4").Value = 4.69
$ws.Range("BI4").Value = 4.97
$ws.Range("BJ4").Value = 5.35
$ws.Range("BK4").Value = 5.14
$ws.Range("BL4").Value = 4.57
$ws.Range("BM4").Value = 5.37
$ws.Range("BN4").Value = 5.27
$ws.Range("BO4").Value = 6.33
$ws.Range("BP4").Value = 5.24
$ws.Range("BQ4").Value = 5.52
$ws.Range("BR4").Value = 5
$ws.Range("BS4").Value = 4.89
$ws.Range("BT4").Value = 5.5
$ws.Range("BU4").Value = 5.54
$ws.Range("BV4").Value = 4.97
$ws.Range("BW4").Value = 5.2
$ws.Range("BX4").Value = 6.26
$ws.Range("BY4").Value = 5.55
$ws.Range("BZ4").Value = 4.79
$ws.Range("CA4").Value = 5.04
$ws.Range("CB4").Value = 5.39
$ws.Range("CC4").Value = 4.53
$ws.Range("CD4").Value = 4.7
$ws.Range("CE4").Value = 4.98
$ws.Range("CF4").Value = 4.63
$ws.Range("CG4").Value = 4.25
$ws.Range("CH4").Value = 6.33
$ws.Range("CI4").Value = 5.23
$ws.Range("CJ4").Value = 6.22
$ws.Range("CK4").Value = 5.34
$ws.Range("CL4").Value = 4.72
$ws.Range("CM4").Value = 4.84
$ws.Range("CN4").Value = 4.75
$ws.Range("CO4").Value = 4.56
$ws.Range("CP4").Value = 5.44
$ws.Range("CQ4").Value = 5.44
$ws.Range("CR4").Value = 5.26
$ws.Range("CS4").Value = 5.26
